# Update odds values for row 34 (match: Oakland Roots - Phoenix Rising)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G34").Value  = 2.32
$ws.Range("H34").Value  = 3.2
$ws.Range("I34").Value  = 2.8
$ws.Range("J34").Value  = 2.92
$ws.Range("L34").Value  = 3.4
$ws.Range("N34").Value  = 7
$ws.Range("U34").Value  = 1.78
$ws.Range("V34").Value  = 1.93
$ws.Range("W34").Value  = 7.7
$ws.Range("X34").Value  = 11.25
$ws.Range("Z34").Value  = 24
$ws.Range("AA34").Value = 19.5
$ws.Range("AC34").Value = 7
$ws.Range("AD34").Value = 6.3
$ws.Range("AG34").Value = 8.75
$ws.Range("AH34").Value = 14
$ws.Range("AI34").Value = 10.25
$ws.Range("AJ34").Value = 35
$ws.Range("AL34").Value = 35
$ws.Range("AN34").Value = 4.3
$ws.Range("AO34").Value = 12.5
$ws.Range("AQ34").Value = 50
$ws.Range("AU34").Value = 7
$ws.Range("AW34").Value = 4.8
$ws.Range("AX34").Value = 15.5
$ws.Range("AY34").Value = 23
$ws.Range("AZ34").Value = 70
